# biotek.xlsx update: rename "ID" sheet to "Channel Map", rework the
# Samples/Channel-Map/Groups/Views header rows and re-point the active
# sheet/selection state to match the author's re-save of the workbook.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Samples": reorder the header/data columns so the Type/Data
# Location/Channels/Plate brand block leads, followed by Plate/Well/Name.
# ---------------------------------------------------------------------
$wsSamples = $wb.Worksheets.Item("Samples")
$wsSamples.Range("A1:G2").Clear()

$wsSamples.Range("A1").Value = "Type"
$wsSamples.Range("B1").Value = "Data Location"
$wsSamples.Range("C1").Value = "Channels"
$wsSamples.Range("D1").Value = "Plate brand"
$wsSamples.Range("E1").Value = "Plate"
$wsSamples.Range("F1").Value = "Well"
$wsSamples.Range("G1").Value = "Name"

$wsSamples.Range("A2").Value = "Plate reader"
$wsSamples.Range("B2").Value = '$GITHUB_WORKSPACE/test/inputs/biotek-data.csv'
$wsSamples.Range("C2").Value = "600,700,(485,530),(485,530[2])"
$wsSamples.Range("C2").NumberFormat = "#,##0"
$wsSamples.Range("D2").Value = "biotek"
$wsSamples.Range("E2").Value = 1

$wsSamples.Columns.Item(1).ColumnWidth = 10.33
$wsSamples.Columns.Item(2).ColumnWidth = 75.5
$wsSamples.Columns.Item(4).ColumnWidth = 13.67

[void]$wsSamples.Range("B6").Select()

# ---------------------------------------------------------------------
# Sheet "ID" -> "Channel Map": header becomes Channel/New name (content
# rows already hold the 600/od1, 700/od2 mapping, unchanged).
# ---------------------------------------------------------------------
$wsChannel = $wb.Worksheets.Item("ID")
$wsChannel.Name = "Channel Map"
$wsChannel.Range("A1").Value = "Channel"
$wsChannel.Range("B1").Value = "New name"
[void]$wsChannel.Range("C3").Select()

# ---------------------------------------------------------------------
# Sheet "Groups": header Group/Name -> Name/Samples.
# ---------------------------------------------------------------------
$wsGroups = $wb.Worksheets.Item("Groups")
$wsGroups.Range("A1").Value = "Name"
$wsGroups.Range("B1").Value = "Samples"
[void]$wsGroups.Range("B2").Select()

# ---------------------------------------------------------------------
# Sheet "Views": header Name/Groups -> Name/View.
# ---------------------------------------------------------------------
$wsViews = $wb.Worksheets.Item("Views")
$wsViews.Range("B1").Value = "View"
[void]$wsViews.Range("B2").Select()
[void]$wsViews.Activate()

Write-Output "done"
